$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ConsumptionSingleValues")

$data = @(
    @("POP_EU_ADULT_2022", "F01", "MeanConsumption", 1.6, "gPerKgBWPerDay", "EU 2022 synthetic"),
    @("POP_EU_ADULT_2022", "F02", "MeanConsumption", 1,   "gPerKgBWPerDay", "EU 2022 synthetic"),
    @("POP_EU_ADULT_2022", "F03", "MeanConsumption", 0.3, "gPerKgBWPerDay", "EU 2022 synthetic"),
    @("POP_EU_ADULT_2022", "F04", "MeanConsumption", 0.8, "gPerKgBWPerDay", "EU 2022 synthetic"),
    @("POP_EU_ADULT_2022", "F05", "MeanConsumption", 1.4, "gPerKgBWPerDay", "EU 2022 synthetic"),
    @("POP_EU_ADULT_2022", "F06", "MeanConsumption", 0.4, "gPerKgBWPerDay", "EU 2022 synthetic"),
    @("POP_EU_ADULT_2022", "F07", "MeanConsumption", 0.3, "gPerKgBWPerDay", "EU 2022 synthetic"),
    @("POP_EU_ADULT_2022", "F08", "MeanConsumption", 2.2, "gPerKgBWPerDay", "EU 2022 synthetic"),
    @("POP_EU_ADULT_2022", "F09", "MeanConsumption", 2.5, "gPerKgBWPerDay", "EU 2022 synthetic"),
    @("POP_EU_ADULT_2022", "F10", "MeanConsumption", 1.8, "gPerKgBWPerDay", "EU 2022 synthetic"),
    @("POP_EU_ADULT_2022", "F11", "MeanConsumption", 0.6, "gPerKgBWPerDay", "EU 2022 synthetic"),
    @("POP_EU_ADULT_2022", "F12", "MeanConsumption", 3,   "gPerKgBWPerDay", "EU 2022 synthetic"),
    @("POP_EU_ADULT_2022", "F13", "MeanConsumption", 1.2, "gPerKgBWPerDay", "EU 2022 synthetic"),
    @("POP_EU_ADULT_2022", "F14", "MeanConsumption", 1.5, "gPerKgBWPerDay", "EU 2022 synthetic"),
    @("POP_EU_ADULT_2022", "F15", "MeanConsumption", 1.3, "gPerKgBWPerDay", "EU 2022 synthetic"),
    @("POP_EU_ADULT_2022", "F16", "MeanConsumption", 0.4, "gPerKgBWPerDay", "EU 2022 synthetic"),
    @("POP_EU_ADULT_2022", "F17", "MeanConsumption", 0.5, "gPerKgBWPerDay", "EU 2022 synthetic"),
    @("POP_EU_ADULT_2022", "F18", "MeanConsumption", 0.2, "gPerKgBWPerDay", "EU 2022 synthetic"),
    @("POP_EU_ADULT_2022", "F19", "MeanConsumption", 1.2, "gPerKgBWPerDay", "EU 2022 synthetic"),
    @("POP_EU_ADULT_2022", "F20", "MeanConsumption", 0.9, "gPerKgBWPerDay", "EU 2022 synthetic")
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = ""
    $ws.Cells.Item($row, 5).Value = $item[3]
    $ws.Cells.Item($row, 6).Value = $item[4]
    $ws.Cells.Item($row, 7).Value = $item[5]
    $row = $row + 1
}
